$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notendasögur")

# Burndown actuals (column I) for the first days of the sprint changed —
# update the "spent today" inputs; the H-column "time left" formulas
# (H3:H16, H{n}=H{n-1}-I{n}) recalculate automatically.
$ws.Range("I3").Value = 80
$ws.Range("I5").Value = 180
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 360
$ws.Range("I8").Value = 420
$ws.Range("I9").Value = 480

# Move the view: scroll back to the top and select F23 instead of I19.
$ws.Activate()
[void]$ws.Range("F23").Select()
